$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("D12").Value = "['MEC-2A-Des. Maq. Cad_T1', -]"

# Row 14
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "['MEC-2A-Des. Maq. Cad_T1', -]"

# Row 15
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "['MEC-2A-Des. Maq. Cad_T1', -]"

# Row 16
$ws.Range("C16").Value = "-"

# Row 18
$ws.Range("C18").Value = "['MEC-1NB-Metalografia', -, 'MEC-1NB-Trat. Termicos', -]"
$ws.Range("D18").Value = "MEC-2NB-Des. Maq. Cad"
$ws.Range("E18").Value = "['MEC-1NA-Metalografia', -, -, -]"
$ws.Range("F18").Value = "-"

# Row 19
$ws.Range("C19").Value = "['MEC-1NB-Metalografia', -, -, -]"
$ws.Range("E19").Value = "['MEC-1NA-Metalografia', -, -, -]"
$ws.Range("F19").Value = "-"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "['MEC-1NB-Trat. Termicos', 'MEC-1NB-Metalografia', -, -]"
$ws.Range("D20").Value = "MEC-2NB-Des. Maq. Cad"
$ws.Range("E20").Value = "['MEC-1NA-Metalografia', -, -, -]"
$ws.Range("F20").Value = "['ELM-1NA-Des. Bas. Mec.', 'ELM-1NA-Des. Bas. Mec.']"

# Row 21
$ws.Range("B21").Value = "[-, -, -, 'MEC-1NB-Trat. Termicos']"
$ws.Range("D21").Value = "['MEC-1NB-Metalografia', 'MEC-1NB-Trat. Termicos', -, -]"
$ws.Range("E21").Value = "['MEC-1NA-Metalografia', -, -, -]"
$ws.Range("F21").Value = "['ELM-1NA-Des. Bas. Mec.', 'ELM-1NA-Des. Bas. Mec.']"
